$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Tgfb2-Tgfbr1 results following Dr Hou advice
$data = @{
    2 = @(3, 1, 2.191602, 6.574806000000001, 0.07674610985252207, 0.07674610985252209, 3, 1, 11.19524, 33.58572, 0.07847249539938134, 0.07847249539938135, 24.53551037448, 220.81959337032, 0.006022458752322453, 0.006022458752322455)
    3 = @(3, 1, 2.191602, 6.574806000000001, 0.07674610985252207, 0.07674610985252209, 3, 1, 31.07813833333333, 93.234415, 0.2178407133195748, 0.2178407133195749, 68.11091012761, 612.99819114849, 0.01671842731477586, 0.01671842731477586)
    4 = @(3, 1, 2.191602, 6.574806000000001, 0.07674610985252207, 0.07674610985252209, 3, 1, 91.56894199999999, 274.706826, 0.6418480872068143, 0.6418480872068144, 200.682676425084, 1806.144087825756, 0.04925934380940534, 0.04925934380940535)
    5 = @(3, 1, 2.191602, 6.574806000000001, 0.07674610985252207, 0.07674610985252209, 3, 1, 8.822188333333335, 26.466565, 0.06183870407422939, 0.06183870407422939, 19.33472559571, 174.01253036139, 0.004745879976018413, 0.004745879976018414)
    6 = @(3, 1, 16.59481266666667, 49.78443799999999, 0.5811216251390648, 0.5811216251390647, 3, 1, 11.19524, 33.58572, 0.07847249539938134, 0.07847249539938135, 185.7829105583733, 1672.04619502536, 0.04560206405520627, 0.04560206405520627)
    7 = @(3, 1, 16.59481266666667, 49.78443799999999, 0.5811216251390648, 0.5811216251390647, 3, 1, 31.07813833333333, 93.234415, 0.2178407133195748, 0.2178407133195749, 515.7358836704188, 4641.62295303377, 0.1265919493457244, 0.1265919493457245)
    8 = @(3, 1, 16.59481266666667, 49.78443799999999, 0.5811216251390648, 0.5811216251390647, 3, 1, 91.56894199999999, 274.706826, 0.6418480872068143, 0.6418480872068144, 1519.569438574865, 13676.12494717378, 0.3729918035300241, 0.3729918035300241)
    9 = @(3, 1, 16.59481266666667, 49.78443799999999, 0.5811216251390648, 0.5811216251390647, 3, 1, 8.822188333333335, 26.466565, 0.06183870407422939, 0.06183870407422939, 146.4025627017189, 1317.62306431547, 0.03593580820810989, 0.03593580820810989)
    10 = @(3, 1, 9.770107666666666, 29.310323, 0.342132265008413, 0.342132265008413, 3, 1, 11.19524, 33.58572, 0.07847249539938134, 0.07847249539938135, 109.3787001541733, 984.40830138756, 0.02684797259185261, 0.02684797259185261)
    11 = @(3, 1, 9.770107666666666, 29.310323, 0.342132265008413, 0.342132265008413, 3, 1, 31.07813833333333, 93.234415, 0.2178407133195748, 0.2178407133195749, 303.6367575962272, 2732.730818366045, 0.07453033665907449, 0.07453033665907452)
    12 = @(3, 1, 9.770107666666666, 29.310323, 0.342132265008413, 0.342132265008413, 3, 1, 91.56894199999999, 274.706826, 0.6418480872068143, 0.6418480872068144, 894.6384222627552, 8051.745800364796, 0.2195969398673848, 0.2195969398673848)
    13 = @(3, 1, 9.770107666666666, 29.310323, 0.342132265008413, 0.342132265008413, 3, 1, 8.822188333333335, 26.466565, 0.06183870407422939, 0.06183870407422939, 86.19372987227723, 775.743568850495, 0.02115701589010108, 0.02115701589010108)
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
